# Edit workbook: add "o_20" and "o_20_jumbled" sheets, add a new
# "evaluator_partial_correctness" column to all sheets, tweak prompt/response
# text on the existing "o_10" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "o_10"

# ---------------------------------------------------------------------------
# Text blocks (kept verbatim, including trailing spaces / blank lines).
# ---------------------------------------------------------------------------

$prompt10 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 10 nodes labelled A to J. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the shortest path from node A to node K? Return the sequence of nodes in response.
   A B C D E F G H I J K
 A 0 1 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0
 I 0 0 0 0 0 0 0 1 0 1 0
 J 0 0 0 0 0 0 0 0 1 0 1
 K 0 0 0 0 0 0 0 0 0 1 0

Solution: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K
        
 Given these examples, answer the following quesiton.

what is the shortest path from node A to node J? Return the sequence of nodes in response.

   A B C D E F G H I J
 A 0 1 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0
 G 0 0 0 0 0 1 0 1 0 0
 H 0 0 0 0 0 0 1 0 1 0
 I 0 0 0 0 0 0 0 1 0 1
 J 0 0 0 0 0 0 0 0 1 0
    
"@

$prompt20 = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the shortest path from node A to node X? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T -> U -> V -> W -> X
 Given these examples, answer the following quesiton.
what is the shortest path from node A to node T? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
"@

$solution10 = "A -> B -> C -> D -> E -> F -> G -> H -> I -> J"
$solution20 = "A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T"

$llmResponse10 = "The shortest path from node A to node J is: A -> B -> C -> D -> E -> F -> G -> H -> I -> J."
$llmResponse20 = "The shortest path from node A to node T is: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T"
$llmResponse20Jumbled = "The shortest path from node A to node T is: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T."

$evaluatorResponse = "invalid input"
$partial20 = "20/20"

# ---------------------------------------------------------------------------
# o_10: add the evaluator_partial_correctness column + tweak existing text.
# ---------------------------------------------------------------------------

$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)

$ws1.Range("A2").Value = $prompt10
$ws1.Range("B2").Value = $solution10
$ws1.Range("C2").Value = $llmResponse10
$ws1.Range("D2").Value = $evaluatorResponse
$ws1.Range("E2").Value = "10/10"
$ws1.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------------
# o_20: new sheet, appended after o_10.
# ---------------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "o_20"

$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A1:E1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)

$ws2.Range("A2").Value = $prompt20
$ws2.Range("B2").Value = $solution20
$ws2.Range("C2").Value = $llmResponse20
$ws2.Range("D2").Value = $evaluatorResponse
$ws2.Range("E2").Value = $partial20
$ws2.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------------
# o_20_jumbled: new sheet, appended after o_20.
# ---------------------------------------------------------------------------

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "o_20_jumbled"

$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A1:E1").Copy()
$ws3.Range("A1:E1").PasteSpecial(-4122)

$ws3.Range("A2").Value = $prompt20
$ws3.Range("B2").Value = $solution20
$ws3.Range("C2").Value = $llmResponse20Jumbled
$ws3.Range("D2").Value = $evaluatorResponse
$ws3.Range("E2").Value = $partial20
$ws3.Rows.Item(2).AutoFit()

$excel.CutCopyMode = $false
$ws1.Select()
